# Update crypto price/volume data per the scraper's Thu Sep 7 22:53:52 UTC 2023 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.190.70'
$ws.Range("E2").Value = '  +1.57%  '
$ws.Range("D3").Value = '1.639.75'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("E4").Value = '  -1.11%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '215.75'
$c.ClearFormats()
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("E6").Value = '  -0.14%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.ClearFormats()
$ws.Range("E7").Value = '  -1.05%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.256'
$c.ClearFormats()
$ws.Range("E8").Value = '  -0.90%  '
$ws.Range("E9").Value = '  -0.61%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.85'
$c.ClearFormats()
$ws.Range("E10").Value = '  +1.03%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0787'
$c.ClearFormats()
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").Value = '1.866.81'
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").Value = '1.607.57'
$ws.Range("E14").Value = '  -1.72%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.553'
$c.ClearFormats()
$ws.Range("E15").Value = '  -1.26%  '
$ws.Range("D16").Value = '0.0₃0765'
$ws.Range("E16").Value = '  -0.34%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '63.36'
$c.ClearFormats()
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("D18").Value = '26.157.33'
$ws.Range("E18").Value = '  +1.41%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.ClearFormats()
$ws.Range("E19").Value = '  -1.00%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '4.48'
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.57%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '194.07'
$c.ClearFormats()
$ws.Range("E21").Value = '  -0.16%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '10.01'
$c.ClearFormats()
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("E23").Value = '  +1.61%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.ClearFormats()
$ws.Range("E24").Value = '  -1.03%  '
$ws.Range("E25").Value = '  -2.49%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '141.74'
$c.ClearFormats()
$ws.Range("E26").Value = '  -0.72%  '
$ws.Range("E27").Value = '  +0.81%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("E35").Value = '  +0.14%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.908'
$c.ClearFormats()
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("D37").Value = '1.142.70'
$ws.Range("E37").Value = '  +1.13%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.548'
$c.ClearFormats()
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  -1.63%  '
$ws.Range("E40").Value = '  +0.18%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.991'
$c.ClearFormats()
$ws.Range("E41").Value = '  -0.91%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.56'
$c.ClearFormats()
$ws.Range("E42").Value = '  -0.53%  '
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("E44").Value = '  -1.98%  '
$ws.Range("D45").Value = '1.776.68'
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("E46").Value = '  -1.45%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '55.87'
$c.ClearFormats()
$ws.Range("E47").Value = '  +1.20%  '
$ws.Range("E48").Value = '  +1.88%  '
$ws.Range("E49").Value = '  +5.09%  '
$ws.Range("E50").Value = '  -0.42%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '7.62'
$c.ClearFormats()
$ws.Range("E51").Value = '  +1.17%  '
